$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.729.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.318.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.52%  "

$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.316.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("E11").Value = "  +2.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.404"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.894.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("E14").Value = "  -2.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.726.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000164"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.311.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "428.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "

$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("E21").Value = "  -1.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.33%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.31%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.452.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("E27").Value = "  +0.44%  "

$ws.Range("E28").Value = "  +6.59%  "

$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("E39").Value = "  -1.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.874.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.45%  "

$ws.Range("E44").Value = "  -3.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("E48").Value = "  -2.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "313.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0271"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.00%  "
